$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 66.04000000000001
$ws.Range("E2").Value = 297.94
$ws.Range("F2").Value = 72.86
$ws.Range("G2").Value = 31.86
$ws.Range("H2").Value = 0.53

$ws.Range("D3").Value = 1.44
$ws.Range("E3").Value = 0.16
$ws.Range("F3").Value = 0.1
$ws.Range("G3").Value = 106.35
$ws.Range("H3").Value = 1.15

$ws.Range("D4").Value = 1.48
$ws.Range("E4").Value = 0.16
$ws.Range("F4").Value = 0.1
$ws.Range("G4").Value = 153.58
$ws.Range("H4").Value = 1.21

$ws.Range("D5").Value = 72.31999999999999
$ws.Range("E5").Value = 0.16
$ws.Range("F5").Value = 0.1
$ws.Range("G5").Value = 160.16
$ws.Range("H5").Value = 1.23

$ws.Range("D6").Value = 81.2
$ws.Range("E6").Value = 0.16
$ws.Range("F6").Value = 0.1
$ws.Range("G6").Value = 137.65
$ws.Range("H6").Value = 1.28
